# "Generate Report for Archive"
#
# The localization-status report was regenerated:
#   - the Status cell that used to read "Ready for handoff" now reads
#     "In Translation" on every sheet (Overview!E2/F2, zh-cn!C2, de-de!C2)
#   - because the new status text is shorter, the Status column got
#     narrower on every sheet that shows it (Overview columns E & F,
#     zh-cn column C, de-de column C)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New width for the (now narrower) Status column, expressed as the
# ColumnWidth (character units) that needs to be assigned so the
# column shrinks to fit the shorter "In Translation" text.
$newStatusColumnWidth = 12.5

# --- Overview sheet: Status lives in columns E and F ------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth

# --- zh-cn sheet: Status lives in column C -----------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth

# --- de-de sheet: Status lives in column C -----------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
